$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 86. This shifts the existing rows 86..161 down to 87..162,
# duplicating the old row 86's formatting (date style) into the now-empty new row 86.
$ws.Rows(86).Insert()

# The old row 86's data is now sitting in row 87 (everything shifted down by one).
# Copy that row's values back into the new, blank row 86 so every column except
# the date (D) and volume (J) keeps its original value.
$ws.Range("A87:R87").Copy()
$ws.Range("A86").PasteSpecial()

# Now overwrite the two cells that actually hold new data for this inserted record.
$ws.Cells.Item(86, 4).Value2 = 44789
$ws.Cells.Item(86, 10).Value2 = 3000
